# Applies the "case with 380 kV done" vm_pu.xlsx update.
# Updates Sheet1 data rows 2-25 (bus voltage-magnitude results) for columns
# B-F and I-N: slack-bus voltage setpoint B changes from 1.05 pu to 1.02 pu,
# and the dependent bus voltages (C-F, I-N) are recomputed accordingly.
# Columns A (bus index), G (=1) and H (blank) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.026300580174367
$ws.Cells.Item(2, 4).Value = 1.031716926897907
$ws.Cells.Item(2, 5).Value = 1.029929704889638
$ws.Cells.Item(2, 6).Value = 1.02481474348338
$ws.Cells.Item(2, 9).Value = 1.03423026287964
$ws.Cells.Item(2, 10).Value = 1.031464788390563
$ws.Cells.Item(2, 11).Value = 1.034524239319187
$ws.Cells.Item(2, 12).Value = 1.032742189863812
$ws.Cells.Item(2, 13).Value = 1.02764214054328
$ws.Cells.Item(2, 14).Value = 1.01450641031337
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.027209720616018
$ws.Cells.Item(3, 4).Value = 1.032412598999369
$ws.Cells.Item(3, 5).Value = 1.030786216798687
$ws.Cells.Item(3, 6).Value = 1.026358023263349
$ws.Cells.Item(3, 9).Value = 1.034473475219142
$ws.Cells.Item(3, 10).Value = 1.03201383797656
$ws.Cells.Item(3, 11).Value = 1.035028836042077
$ws.Cells.Item(3, 12).Value = 1.033406820359576
$ws.Cells.Item(3, 13).Value = 1.028990591554344
$ws.Cells.Item(3, 14).Value = 1.014689793002341
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.027797934551339
$ws.Cells.Item(4, 4).Value = 1.032862567600095
$ws.Cells.Item(4, 5).Value = 1.031340727116452
$ws.Cells.Item(4, 6).Value = 1.027356556255189
$ws.Cells.Item(4, 9).Value = 1.034629400442383
$ws.Cells.Item(4, 10).Value = 1.032368436995784
$ws.Cells.Item(4, 11).Value = 1.035354492521838
$ws.Cells.Item(4, 12).Value = 1.033836532024417
$ws.Cells.Item(4, 13).Value = 1.029862577667189
$ws.Cells.Item(4, 14).Value = 1.014808183708705
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028045205094899
$ws.Cells.Item(5, 4).Value = 1.033051691250667
$ws.Cells.Item(5, 5).Value = 1.031573912240471
$ws.Cells.Item(5, 6).Value = 1.027776327152223
$ws.Cells.Item(5, 9).Value = 1.034694604259561
$ws.Cells.Item(5, 10).Value = 1.032517349293383
$ws.Cells.Item(5, 11).Value = 1.03549119445
$ws.Cells.Item(5, 12).Value = 1.034017099127604
$ws.Cells.Item(5, 13).Value = 1.030229032686768
$ws.Cells.Item(5, 14).Value = 1.014857890408499
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028086722049632
$ws.Cells.Item(6, 4).Value = 1.033083443396179
$ws.Cells.Item(6, 5).Value = 1.031613069082482
$ws.Cells.Item(6, 6).Value = 1.027846807950149
$ws.Cells.Item(6, 9).Value = 1.034705531906078
$ws.Cells.Item(6, 10).Value = 1.032542342872658
$ws.Cells.Item(6, 11).Value = 1.035514135332167
$ws.Cells.Item(6, 12).Value = 1.03404741222581
$ws.Cells.Item(6, 13).Value = 1.030290554715623
$ws.Cells.Item(6, 14).Value = 1.01486623258268
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.027801238649869
$ws.Cells.Item(7, 4).Value = 1.032865094849873
$ws.Cells.Item(7, 5).Value = 1.031343842676937
$ws.Cells.Item(7, 6).Value = 1.027362165296957
$ws.Cells.Item(7, 9).Value = 1.034630273063223
$ws.Cells.Item(7, 10).Value = 1.032370427402417
$ws.Cells.Item(7, 11).Value = 1.035356319941809
$ws.Cells.Item(7, 12).Value = 1.033838945099955
$ws.Cells.Item(7, 13).Value = 1.029867474754633
$ws.Cells.Item(7, 14).Value = 1.014808848146597
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.02660784098305
$ws.Cells.Item(8, 4).Value = 1.031952068996278
$ws.Cells.Item(8, 5).Value = 1.030219106846783
$ws.Cells.Item(8, 6).Value = 1.025336320026476
$ws.Cells.Item(8, 9).Value = 1.034312757650815
$ws.Cells.Item(8, 10).Value = 1.031650481533599
$ws.Cells.Item(8, 11).Value = 1.034694946275553
$ws.Cells.Item(8, 12).Value = 1.032966876936018
$ws.Cells.Item(8, 13).Value = 1.028097972656346
$ws.Cells.Item(8, 14).Value = 1.014568441243381
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02450445856448
$ws.Cells.Item(9, 4).Value = 1.03034186478378
$ws.Cells.Item(9, 5).Value = 1.028239415486766
$ws.Cells.Item(9, 6).Value = 1.02176574847337
$ws.Cells.Item(9, 9).Value = 1.033742160415935
$ws.Cells.Item(9, 10).Value = 1.030376694558871
$ws.Cells.Item(9, 11).Value = 1.03352301127829
$ws.Cells.Item(9, 12).Value = 1.031427518869811
$ws.Cells.Item(9, 13).Value = 1.024975470274806
$ws.Cells.Item(9, 14).Value = 1.014142747698579
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.023101891385902
$ws.Cells.Item(10, 4).Value = 1.029267529235544
$ws.Cells.Item(10, 5).Value = 1.026921145546334
$ws.Cells.Item(10, 6).Value = 1.019384538325242
$ws.Cells.Item(10, 9).Value = 1.033354310079525
$ws.Cells.Item(10, 10).Value = 1.029524041164258
$ws.Cells.Item(10, 11).Value = 1.032737355858804
$ws.Cells.Item(10, 12).Value = 1.030399498527726
$ws.Cells.Item(10, 13).Value = 1.022890571668311
$ws.Cells.Item(10, 14).Value = 1.013857566590638
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.022494489052901
$ws.Cells.Item(11, 4).Value = 1.028802130752149
$ws.Cells.Item(11, 5).Value = 1.026350685096326
$ws.Cells.Item(11, 6).Value = 1.018353182203156
$ws.Cells.Item(11, 9).Value = 1.033184600079311
$ws.Cells.Item(11, 10).Value = 1.029154011060142
$ws.Cells.Item(11, 11).Value = 1.03239612473031
$ws.Cells.Item(11, 12).Value = 1.02995393278361
$ws.Cells.Item(11, 13).Value = 1.021986962411386
$ws.Cells.Item(11, 14).Value = 1.013733751923973
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022268860320855
$ws.Cells.Item(12, 4).Value = 1.028629230781779
$ws.Cells.Item(12, 5).Value = 1.026138844826081
$ws.Cells.Item(12, 6).Value = 1.017970043064049
$ws.Cells.Item(12, 9).Value = 1.033121296584022
$ws.Cells.Item(12, 10).Value = 1.029016441229827
$ws.Cells.Item(12, 11).Value = 1.032269220532229
$ws.Cells.Item(12, 12).Value = 1.029788365788502
$ws.Cells.Item(12, 13).Value = 1.021651191188603
$ws.Cells.Item(12, 14).Value = 1.013687712170366
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022317259017242
$ws.Cells.Item(13, 4).Value = 1.028666319774775
$ws.Cells.Item(13, 5).Value = 1.026184282830293
$ws.Cells.Item(13, 6).Value = 1.018052229958919
$ws.Cells.Item(13, 9).Value = 1.033134887419447
$ws.Cells.Item(13, 10).Value = 1.029045956046924
$ws.Cells.Item(13, 11).Value = 1.032296448961303
$ws.Cells.Item(13, 12).Value = 1.029823883369182
$ws.Cells.Item(13, 13).Value = 1.021723221244551
$ws.Cells.Item(13, 14).Value = 1.013697590092449
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.022475838767005
$ws.Cells.Item(14, 4).Value = 1.028787839401436
$ws.Cells.Item(14, 5).Value = 1.026333173208009
$ws.Cells.Item(14, 6).Value = 1.018321512807426
$ws.Cells.Item(14, 9).Value = 1.033179372811647
$ws.Cells.Item(14, 10).Value = 1.029142642029016
$ws.Cells.Item(14, 11).Value = 1.032385637966938
$ws.Cells.Item(14, 12).Value = 1.029940248270627
$ws.Cells.Item(14, 13).Value = 1.021959210149633
$ws.Cells.Item(14, 14).Value = 1.013729947272775
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.022573543326812
$ws.Cells.Item(15, 4).Value = 1.028862707657605
$ws.Cells.Item(15, 5).Value = 1.026424916663282
$ws.Cells.Item(15, 6).Value = 1.018487420375787
$ws.Cells.Item(15, 9).Value = 1.033206746529832
$ws.Cells.Item(15, 10).Value = 1.029202196994968
$ws.Cells.Item(15, 11).Value = 1.032440569618296
$ws.Cells.Item(15, 12).Value = 1.030011936040611
$ws.Cells.Item(15, 13).Value = 1.022104593254164
$ws.Cells.Item(15, 14).Value = 1.013749877044022
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.02314220094293
$ws.Cells.Item(16, 4).Value = 1.029298411935327
$ws.Cells.Item(16, 5).Value = 1.026959012710122
$ws.Cells.Item(16, 6).Value = 1.01945297969681
$ws.Cells.Item(16, 9).Value = 1.033365535919677
$ws.Cells.Item(16, 10).Value = 1.029548581437951
$ws.Cells.Item(16, 11).Value = 1.032759980398403
$ws.Cells.Item(16, 12).Value = 1.030429060279115
$ws.Cells.Item(16, 13).Value = 1.022950523219527
$ws.Cells.Item(16, 14).Value = 1.013865776823838
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.023498882872783
$ws.Cells.Item(17, 4).Value = 1.029571663367648
$ws.Cells.Item(17, 5).Value = 1.027294133333068
$ws.Cells.Item(17, 6).Value = 1.020058571299832
$ws.Cells.Item(17, 9).Value = 1.033464666910215
$ws.Cells.Item(17, 10).Value = 1.02976563812333
$ws.Cells.Item(17, 11).Value = 1.032960060709901
$ws.Cells.Item(17, 12).Value = 1.030690597192778
$ws.Cells.Item(17, 13).Value = 1.023480925714934
$ws.Cells.Item(17, 14).Value = 1.013938389521027
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.023706921560327
$ws.Cells.Item(18, 4).Value = 1.029731026585086
$ws.Cells.Item(18, 5).Value = 1.027489638161862
$ws.Cells.Item(18, 6).Value = 1.020411776371524
$ws.Cells.Item(18, 9).Value = 1.033522317729044
$ws.Cells.Item(18, 10).Value = 1.029892163971221
$ws.Cells.Item(18, 11).Value = 1.03307666409742
$ws.Cells.Item(18, 12).Value = 1.030843106133515
$ws.Cells.Item(18, 13).Value = 1.023790220355712
$ws.Cells.Item(18, 14).Value = 1.013980711499313
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.023777856052041
$ws.Cells.Item(19, 4).Value = 1.029785361966393
$ws.Cells.Item(19, 5).Value = 1.027556306125847
$ws.Cells.Item(19, 6).Value = 1.020532205901381
$ws.Cells.Item(19, 9).Value = 1.033541946229415
$ws.Cells.Item(19, 10).Value = 1.02993529251083
$ws.Cells.Item(19, 11).Value = 1.033116405863332
$ws.Cells.Item(19, 12).Value = 1.030895100767443
$ws.Cells.Item(19, 13).Value = 1.023895668395053
$ws.Cells.Item(19, 14).Value = 1.013995136809525
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.023460615090451
$ws.Cells.Item(20, 4).Value = 1.029542348111634
$ws.Cells.Item(20, 5).Value = 1.027258174477418
$ws.Cells.Item(20, 6).Value = 1.019993599811629
$ws.Cells.Item(20, 9).Value = 1.033454048740002
$ws.Cells.Item(20, 10).Value = 1.029742358234106
$ws.Cells.Item(20, 11).Value = 1.032938604336539
$ws.Cells.Item(20, 12).Value = 1.030662541009923
$ws.Cells.Item(20, 13).Value = 1.023424026874932
$ws.Cells.Item(20, 14).Value = 1.01393060215281
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.022429141328959
$ws.Cells.Item(21, 4).Value = 1.02875205573201
$ws.Cells.Item(21, 5).Value = 1.026289327206773
$ws.Cells.Item(21, 6).Value = 1.018242217081481
$ws.Cells.Item(21, 9).Value = 1.033166280305414
$ws.Cells.Item(21, 10).Value = 1.029114173846117
$ws.Cells.Item(21, 11).Value = 1.032359378327175
$ws.Cells.Item(21, 12).Value = 1.029905983475407
$ws.Cells.Item(21, 13).Value = 1.021889720943189
$ws.Cells.Item(21, 14).Value = 1.013720420255295
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.021780540628772
$ws.Cells.Item(22, 4).Value = 1.028254992716773
$ws.Cells.Item(22, 5).Value = 1.025680487801258
$ws.Cells.Item(22, 6).Value = 1.017140775111957
$ws.Cells.Item(22, 9).Value = 1.032983811934649
$ws.Cells.Item(22, 10).Value = 1.028718490880772
$ws.Cells.Item(22, 11).Value = 1.031994294768585
$ws.Cells.Item(22, 12).Value = 1.029429935204327
$ws.Cells.Item(22, 13).Value = 1.020924284242312
$ws.Cells.Item(22, 14).Value = 1.013587984318667
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.022124383050627
$ws.Cells.Item(23, 4).Value = 1.028518511662381
$ws.Cells.Item(23, 5).Value = 1.026003215330325
$ws.Cells.Item(23, 6).Value = 1.017724698722992
$ws.Cells.Item(23, 9).Value = 1.033080687550158
$ws.Cells.Item(23, 10).Value = 1.028928318063028
$ws.Cells.Item(23, 11).Value = 1.032187917874738
$ws.Cells.Item(23, 12).Value = 1.029682332527736
$ws.Cells.Item(23, 13).Value = 1.021436154048236
$ws.Cells.Item(23, 14).Value = 1.013658218243617
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.023477906678137
$ws.Cells.Item(24, 4).Value = 1.029555594474846
$ws.Cells.Item(24, 5).Value = 1.027274422629325
$ws.Cells.Item(24, 6).Value = 1.020022957710225
$ws.Cells.Item(24, 9).Value = 1.033458847161153
$ws.Cells.Item(24, 10).Value = 1.029752877659926
$ws.Cells.Item(24, 11).Value = 1.032948299857192
$ws.Cells.Item(24, 12).Value = 1.030675218518209
$ws.Cells.Item(24, 13).Value = 1.023449737258
$ws.Cells.Item(24, 14).Value = 1.013934121027461
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025048288293523
$ws.Cells.Item(25, 4).Value = 1.030758296145212
$ws.Cells.Item(25, 5).Value = 1.028750945949641
$ws.Cells.Item(25, 6).Value = 1.022688950223602
$ws.Cells.Item(25, 9).Value = 1.033890987487564
$ws.Cells.Item(25, 10).Value = 1.030706609271034
$ws.Cells.Item(25, 11).Value = 1.033826754827773
$ws.Cells.Item(25, 12).Value = 1.031825794451045
$ws.Cells.Item(25, 13).Value = 1.025783262036998
$ws.Cells.Item(25, 14).Value = 1.014253044124094
